# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q1"),
#    carrying the same look/formatting as the other quarterly fund sheets.
# 2) Update the "总计" summary sheet with a new row for 2022-Q4 and shift the
#    existing rows down, renumbering the serial index column.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# --- 1) Create the new "2022-Q4" sheet -------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"

# Copy the header row + index-column formatting from the existing "2022-Q1"
# sheet so the new sheet matches the workbook's existing look.
$q1Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1Sheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move it into the correct tab position: right after "总计", before "2022-Q1"
$newSheet.Move($q1Sheet)

# NOTE: after .Move() the old $newSheet reference becomes stale (it starts
# resolving to whatever sheet now sits at the old index), so re-fetch it by
# its (now unique) name before doing any further work on it.
$newSheet = $wb.Worksheets.Item("2022-Q4")

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2: 161039 fund data
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "161039"
$newSheet.Range("C2").Value = "富国中证1000指数增强（LOF）A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "20.17"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "91.85"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "0.73"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.1472"
$newSheet.Range("H2").Value = 3

# Row 3: 013331 fund data
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "013331"
$newSheet.Range("C3").Value = "富国中证1000指数增强（LOF）C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "6.90"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "91.85"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "0.73"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0504"
$newSheet.Range("H3").Value = 3

# --- 2) Update the "总计" summary sheet -------------------------------------
# Shift the existing 3 data rows (currently rows 2-4) down to rows 3-5 by
# rewriting their contents directly (avoids unpredictable formatting from a
# row-insert), then fill row 2 with the new "2022-Q4" entry. Finally,
# (re)apply column A's existing index-cell formatting across A2:A5 so every
# row in the serial-number column looks the same, matching the source sheet.

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q1"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.13

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.06

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.06

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.2

# Row 5 is brand new; give its index cell (A5) the same look as the other
# index cells in column A (bold / bordered / centered).
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
